$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-160 down to 51-161
$ws.Rows("50:50").Insert()

# Populate the newly inserted row 50 with the new record
$ws.Cells.Item(50, 1).Value = 10
$ws.Cells.Item(50, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(50, 3).Value = "La Araucanía"
$ws.Cells.Item(50, 4).Value = 44498
$ws.Cells.Item(50, 5).Value = 9
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100102
$ws.Cells.Item(50, 8).Value = "Cítricos"
$ws.Cells.Item(50, 9).Value = 100102006
$ws.Cells.Item(50, 10).Value = "Pomelo"
$ws.Cells.Item(50, 11).Value = "Red Blush"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 65
$ws.Cells.Item(50, 14).Value = 12000
$ws.Cells.Item(50, 15).Value = 12000
$ws.Cells.Item(50, 16).Value = 12000
$ws.Cells.Item(50, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(50, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(50, 19).Value = 800
$ws.Cells.Item(50, 20).Value = 15
